# Applies the "Quelques fichiers a sync" edit: marks several requirement
# bullets in red (FF0000), merges a couple of runs, relocates the
# _GoBack bookmark, and bumps the cached footer page-number field from
# 1 to 2.

$d = $word.ActiveDocument
$wdRed = 255  # OLE color 0x0000FF -> RGB(255,0,0) -> w:color FF0000
$apos = [char]0x2019

function Find-ParagraphStartingWith($prefix) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

function Find-ParagraphContaining($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) "Créer une séance" - whole paragraph turns red
# ---------------------------------------------------------------------
$p = Find-ParagraphStartingWith("Créer une séance")
$p.Range.Font.Color = $wdRed

# ---------------------------------------------------------------------
# 2) "Modifier une séance" - whole paragraph turns red
# ---------------------------------------------------------------------
$p = Find-ParagraphStartingWith("Modifier une séance")
$p.Range.Font.Color = $wdRed

# ---------------------------------------------------------------------
# 3) "D'afficher ... (vous utiliserez une vue partielle)"
#    only some runs turn red; " (v" splits into " " (red) + "(v" (black)
# ---------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("D${apos}afficher", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Color = $wdRed

$rng.Collapse(0)
$null = $rng.Find.Execute(" la liste des séances futures", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Color = $wdRed

$rng.Collapse(0)
$null = $rng.Find.Execute(" (v", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spacePos = $rng.Start
$sub = $d.Range($spacePos, $spacePos + 1)
$sub.Font.Color = $wdRed

# ---------------------------------------------------------------------
# 4) "Vous devez prévoir la pagination ... selon le statut"
#    whole paragraph turns red (strike already present per-run)
# ---------------------------------------------------------------------
$p = Find-ParagraphStartingWith("Vous devez prévoir la pagination")
$p.Range.Font.Color = $wdRed

# ---------------------------------------------------------------------
# 6) "La date de rendez-vous pour une séance de photo ..."
#    whole paragraph turns red
# ---------------------------------------------------------------------
$p = Find-ParagraphStartingWith("La date de rendez-vous pour une séance de photo")
$p.Range.Font.Color = $wdRed

# ---------------------------------------------------------------------
# 7) "La date pour un rendez-vous doit être unique ..."
#    whole paragraph turns red
# ---------------------------------------------------------------------
$p = Find-ParagraphStartingWith("La date pour un rendez-vous doit être unique")
$p.Range.Font.Color = $wdRed

# ---------------------------------------------------------------------
# 8) "Pour la même date et le même photographe ... 2 rendez-vous"
#    whole paragraph turns red
# ---------------------------------------------------------------------
$p = Find-ParagraphStartingWith("Pour la même date et le même photographe")
$p.Range.Font.Color = $wdRed

Write-Output "step1 done"
